$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix wrong function-key assignment: strip the stray "c" from Fc1, Fc2, Fc6, Fc7, Fc9
$ws.Range("A3").Value = "F1"
$ws.Range("A4").Value = "F2"
$ws.Range("A5").Value = "F6"
$ws.Range("A6").Value = "F7"
$ws.Range("A7").Value = "F9"
